$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Timestamp update ---
$ws.Range("A1").Value = "Datos actualizados a 5 de Julio de 2020 a las 12:53"

# --- Row 4 : Estados Unidos (data refresh) ---
$ws.Range("B4").Value = 2936122
$ws.Range("C4").Value = 352
$ws.Range("E4").Value = 1543185

# --- Row 7 : India (data refresh) ---
$ws.Range("B7").Value = 675453
$ws.Range("C7").Value = 1549
$ws.Range("E7").Value = 247067
$ws.Range("G7").Value = 24
$ws.Range("H7").Value = 19303

# --- Row 14 : Iran (data refresh) ---
$ws.Range("B14").Value = 240438
$ws.Range("C14").Value = 2560
$ws.Range("D14").Value = 201330
$ws.Range("E14").Value = 27537
$ws.Range("G14").Value = 163
$ws.Range("H14").Value = 11571

# --- Row 21 : Banglades (data refresh) ---
$ws.Range("B21").Value = 162417
$ws.Range("C21").Value = 2738
$ws.Range("D21").Value = 72625
$ws.Range("E21").Value = 87740
$ws.Range("G21").Value = 55
$ws.Range("H21").Value = 2052

# --- Rows 29/30 : Indonesia overtakes Bielorrusia ---
$ws.Range("A29").Value = "Indonesia"
$ws.Range("B29").Value = 63749
$ws.Range("C29").Value = 1607
$ws.Range("D29").Value = 29105
$ws.Range("E29").Value = 31473
$ws.Range("G29").Value = 82
$ws.Range("H29").Value = 3171

$ws.Range("A30").Value = "Bielorrusia"
$ws.Range("B30").Value = 63270
$ws.Range("D30").Value = 50669
$ws.Range("E30").Value = 12183
$ws.Range("H30").Value = 418

# --- Row 37 : Ucrania (data refresh) ---
$ws.Range("B37").Value = 48500
$ws.Range("C37").Value = 823
$ws.Range("D37").Value = 21376
$ws.Range("E37").Value = 25875
$ws.Range("G37").Value = 22
$ws.Range("H37").Value = 1249

# --- Row 39 : Oman (data refresh) ---
$ws.Range("B39").Value = 46178
$ws.Range("C39").Value = 1072
$ws.Range("D39").Value = 27917
$ws.Range("E39").Value = 18048
$ws.Range("G39").Value = 10
$ws.Range("H39").Value = 213

# --- Rows 41/42 : Filipinas overtakes Portugal ---
$ws.Range("A41").Value = "Filipinas"
$ws.Range("B41").Value = 44254
$ws.Range("C41").Value = 2424
$ws.Range("D41").Value = 11942
$ws.Range("E41").Value = 31015
$ws.Range("G41").Value = 7
$ws.Range("H41").Value = 1297

$ws.Range("A42").Value = "Portugal"
$ws.Range("B42").Value = 43569
$ws.Range("D42").Value = 28772
$ws.Range("E42").Value = 13192
$ws.Range("H42").Value = 1605

# --- Rows 50/51/52 : Rumania overtakes Barein and Armenia ---
$ws.Range("A50").Value = "Rumania"
$ws.Range("B50").Value = 28973
$ws.Range("C50").Value = 391
$ws.Range("D50").Value = 20026
$ws.Range("E50").Value = 7197
$ws.Range("G50").Value = 19
$ws.Range("H50").Value = 1750

$ws.Range("A51").Value = "Barein"
$ws.Range("B51").Value = 28857
$ws.Range("C51").Value = 0
$ws.Range("D51").Value = 23959
$ws.Range("E51").Value = 4802
$ws.Range("G51").Value = 0
$ws.Range("H51").Value = 96

$ws.Range("A52").Value = "Armenia"
$ws.Range("B52").Value = 28606
$ws.Range("C52").Value = 706
$ws.Range("D52").Value = 16140
$ws.Range("E52").Value = 11982
$ws.Range("G52").Value = 7
$ws.Range("H52").Value = 484

# --- Rows 57/58 : Ghana overtakes Azerbaiyan ---
$ws.Range("A57").Value = "Ghana"
$ws.Range("B57").Value = 20085
$ws.Range("C57").Value = 697
$ws.Range("D57").Value = 14870
$ws.Range("E57").Value = 5093
$ws.Range("G57").Value = 5
$ws.Range("H57").Value = 122

$ws.Range("A58").Value = "Azerbaiyan"
$ws.Range("B58").Value = 19801
$ws.Range("D58").Value = 11291
$ws.Range("E58").Value = 8269
$ws.Range("H58").Value = 241

# --- Row 73 : Noruega (data refresh) ---
$ws.Range("B73").Value = 8927
$ws.Range("C73").Value = 1
$ws.Range("E73").Value = 538

# --- Row 114 : Estonia (data refresh) ---
$ws.Range("D114").Value = 1874
$ws.Range("E114").Value = 50

# --- Rows 205/206 : Fiyi overtakes Dominica (tied totals, only order changes) ---
$ws.Range("A205").Value = "Fiyi"
$ws.Range("A206").Value = "Dominica"
